# ENH: finish work of LexicalAnalyzer
#
# Adds the missing "空格" (space) transition column (CB) to the DFA table
# on Sheet1. Column CB mirrors column CA (the default/"any other char"
# column) for every state row, except:
#   - CB1 is the column header "空格" (instead of "}")
#   - CB2 (state 1, on space) loops back to state 2, instead of going to
#     state 25 like CA2 does.
# Rows 20-26 already contained empty (but styled) CB/CC/CD placeholder
# cells, so those just get their value filled in; rows 1-19 need the CB
# cell created from scratch.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (row 1) ------------------------------------------------------
$ws.Cells.Item(1, 80).Value2 = "空格"

# --- Numeric transition rows (2-19) --------------------------------------
# These are typed as real numbers first (so the stored cell value is
# numeric), then the column is formatted as Text ("@") to match the rest
# of the table's look - mirroring how the other transition columns are
# built (value first, text format applied after, so the stored <v> stays
# numeric even though the cell displays/format as text).
$spaceCol = @{
    2  = 1
    3  = -1
    4  = -2
    5  = -3
    6  = -4
    7  = -5
    8  = -7
    9  = 8
    10 = 0
    11 = -9
    12 = -8
    13 = -11
    14 = -10
    15 = -13
    16 = -12
    17 = 0
    18 = -14
    19 = -15
}

foreach ($row in $spaceCol.Keys) {
    $cell = $ws.Cells.Item($row, 80)
    $cell.Value2 = $spaceCol[$row]
    $cell.NumberFormat = "@"
}

# --- Accepting-state rows (20-26) -----------------------------------------
# CB20:CB26 already exist as blank styled cells in the sheet; just give
# them the same text values already shown in column CA.
$acceptCol = @{
    20 = "-16"
    21 = "-17"
    22 = "-18"
    23 = "-19"
    24 = "-20"
    25 = "-21"
    26 = "-22"
}

foreach ($row in $acceptCol.Keys) {
    $ws.Cells.Item($row, 80).Value2 = $acceptCol[$row]
}

# --- View / selection ------------------------------------------------------
# Scroll so the new column is visible and select the (now one column wider)
# table range, matching the saved view state of the finished sheet.
try {
    $excel.ActiveWindow.ScrollColumn = 64
} catch {
}

$selectResult = $ws.Range("B2:CB26").Select()
